$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in cell E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active cell/selection
$ws.Range("E8").Select()
